$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy the style of the existing header cell (H1) so the
# new headers match formatting (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-31
$data = @(
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(7, 8),
    @(8, 9),
    @(6, 7),
    @(7, 8),
    @(7, 8),
    @(10, 10),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(4, 4),
    @(6, 7),
    @(5, 6),
    @(6, 7),
    @(6, 6),
    @(7, 8),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(2, 3),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$wb.Save()
